$d = $word.ActiveDocument

# 1. Title paragraph: merge the word-by-word runs back into a single run.
$titleRange = $d.Paragraphs(1).Range
$titleRange.Find.Execute(
    "Answers: Using the quadratic formula", $true, $false, $false, $false, $false,
    $true, 1, $false, "Answers: Using the quadratic formula", 2)

# 2. Author paragraph: merge "Tom" / " " / "Coleman" into a single run.
$authorRange = $d.Paragraphs(2).Range
$authorRange.Find.Execute(
    "Tom Coleman", $true, $false, $false, $false, $false,
    $true, 1, $false, "Tom Coleman", 2)

# 3. Abstract paragraph: merge the word-by-word runs back into a single run.
$abstractRange = $d.Paragraphs(4).Range
$abstractRange.Find.Execute(
    "Answers to questions relating to the guide on using the quadratic formula.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Answers to questions relating to the guide on using the quadratic formula.", 2)
